$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = 28
$ws.Range("B12").Value = "2：44-6;57"
$ws.Range("C12").Value = "关键字 ifelse switch while for（未结束）"

$ws.Range("C12").Select()
